$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New credit/attribution row for the projectile-hit sound effect
$ws.Range("A2").Value = "87535__whiprealgood__splat"
$ws.Range("B2").Value = "https://freesound.org/people/Whiprealgood/sounds/87535/"
$ws.Range("C2").Value = "Creative Commons License"
$ws.Range("D2").Value = "Doesn't require attribution"

# Scroll the view over a bit and land the selection on the new row's
# last cell, matching where the author's cursor ended up after entry.
try { $excel.ActiveWindow.ScrollColumn = 2 } catch { }
$ws.Range("D2").Select()
